$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.039.49"
$ws.Range("E2").Value = "  -4.23%  "
$ws.Range("D3").Value = "3.234.34"
$ws.Range("E3").Value = "  -7.68%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'593.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "'153.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.12%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.228.11"
$ws.Range("E8").Value = "  -7.74%  "
$ws.Range("E9").Value = "  -10.02%  "
$ws.Range("D10").Value = "'0.175"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.68%  "
$ws.Range("D11").Value = "'6.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.45%  "
$ws.Range("D12").Value = "'0.505"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -13.73%  "
$ws.Range("D13").Value = "'39.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -14.32%  "
$ws.Range("D14").Value = "'0.0000248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.69%  "
$ws.Range("D15").Value = "3.751.39"
$ws.Range("E15").Value = "  -7.82%  "
$ws.Range("D16").Value = "67.070.72"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").Value = "3.230.70"
$ws.Range("E17").Value = "  -7.74%  "
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "'7.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -12.95%  "
$ws.Range("D20").Value = "'535.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.74%  "
$ws.Range("D21").Value = "'15.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -13.51%  "
$ws.Range("D22").Value = "'0.766"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.68%  "
$ws.Range("D23").Value = "'7.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.69%  "
$ws.Range("D24").Value = "'13.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.26%  "
$ws.Range("D25").Value = "'86.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.50%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'3.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -14.03%  "
$ws.Range("E28").Value = "  -13.12%  "
$ws.Range("D29").Value = "'8.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.09%  "
$ws.Range("D30").Value = "'29.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.33%  "
$ws.Range("D31").Value = "'2.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.98%  "
$ws.Range("E32").Value = "  -10.14%  "
$ws.Range("D33").Value = "'541.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -14.38%  "
$ws.Range("D34").Value = "'6.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -17.75%  "
$ws.Range("D35").Value = "'5.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -14.63%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("D38").Value = "'0.0877"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -12.02%  "
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").Value = "'9.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.32%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0426"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.26%  "
$ws.Range("D41").Value = "'0.128"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.46%  "
$ws.Range("D42").Value = "2.960.21"
$ws.Range("E42").Value = "  -11.73%  "
$ws.Range("D43").Value = "'2.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -22.95%  "
$ws.Range("D44").Value = "'0.269"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -13.29%  "
$ws.Range("D45").Value = "0.0₃0595"
$ws.Range("E45").Value = "  -18.78%  "
$ws.Range("D46").Value = "'2.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -15.97%  "
$ws.Range("D47").Value = "'26.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -16.26%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'2.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -16.14%  "
$ws.Range("D51").Value = "'122.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.07%  "
